$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-RowRange {
    param($ws, $rowNum, $startCol, $vals)
    $n = $vals.Count
    $arr = New-Object 'object[,]' 1,$n
    for ($i = 0; $i -lt $n; $i++) {
        $arr[0,$i] = $vals[$i]
    }
    $startColIdx = $ws.Range($startCol + "1").Column
    $endColIdx = $startColIdx + $n - 1
    $startCell = $ws.Cells.Item($rowNum, $startColIdx)
    $endCell = $ws.Cells.Item($rowNum, $endColIdx)
    $rng = $ws.Range($startCell, $endCell)
    $rng.Value = $arr
}

Set-RowRange $ws 2 "G" @(0.3201913333333333,0.960574,0.03939146858412543,0.03939146858412543)
Set-RowRange $ws 2 "M" @(8.820647333333334,26.461942,0.06415146660411865,0.06415146660411865,2.824294830523111,25.418653474708,0.002527020481361712,0.002527020481361711)
Set-RowRange $ws 3 "G" @(0.3201913333333333,0.960574,0.03939146858412543,0.03939146858412543)
Set-RowRange $ws 3 "O" @(0.3979101621202897,0.3979101621202898,17.51815934036333,157.66343406327,0.01567426565046565,0.01567426565046565)
Set-RowRange $ws 4 "G" @(0.3201913333333333,0.960574,0.03939146858412543,0.03939146858412543)
Set-RowRange $ws 4 "M" @(21.90816333333333,65.72449,0.1593353362087987,0.1593353362087987,7.014804028584445,63.13323625726,0.006276452890609957,0.006276452890609956)
Set-RowRange $ws 5 "G" @(0.3201913333333333,0.960574,0.03939146858412543,0.03939146858412543)
Set-RowRange $ws 5 "M" @(52.056859,156.170577,0.3786030350667928,0.3786030350667929,16.66815509235533,150.013395831198,0.01491372956168811,0.01491372956168811)
Set-RowRange $ws 6 "I" @(0.2346323697636092,0.2346323697636091)
Set-RowRange $ws 6 "M" @(8.820647333333334,26.461942,0.06415146660411865,0.06415146660411865,16.822703311544,151.404329803896,0.01505201063313539,0.01505201063313539)
Set-RowRange $ws 7 "I" @(0.2346323697636092,0.2346323697636091)
Set-RowRange $ws 7 "O" @(0.3979101621202897,0.3979101621202898)
Set-RowRange $ws 7 "S" @(0.09336260429130548,0.09336260429130548)
Set-RowRange $ws 8 "I" @(0.2346323697636092,0.2346323697636091)
Set-RowRange $ws 8 "M" @(21.90816333333333,65.72449,0.1593353362087987,0.1593353362087987,41.78316147668001,376.0484532901201,0.03738522752175183,0.03738522752175182)
Set-RowRange $ws 9 "I" @(0.2346323697636092,0.2346323697636091)
Set-RowRange $ws 9 "M" @(52.056859,156.170577,0.3786030350667928,0.3786030350667929,99.282633257364,893.5436993162759,0.08883252731741642,0.08883252731741641)
Set-RowRange $ws 10 "G" @(5.780535,17.341605,0.7111490510422025,0.7111490510422023)
Set-RowRange $ws 10 "M" @(8.820647333333334,26.461942,0.06415146660411865,0.06415146660411865,50.98806063299001,458.89254569691,0.04562125459848452,0.0456212545984845)
Set-RowRange $ws 11 "G" @(5.780535,17.341605,0.7111490510422025,0.7111490510422023)
Set-RowRange $ws 11 "O" @(0.3979101621202897,0.3979101621202898,316.261942971225,2846.357486741025,0.282973434191893,0.282973434191893)
Set-RowRange $ws 12 "G" @(5.780535,17.341605,0.7111490510422025,0.7111490510422023)
Set-RowRange $ws 12 "M" @(21.90816333333333,65.72449,0.1593353362087987,0.1593353362087987,126.64090493405,1139.76814440645,0.1133111731423775,0.1133111731423774)
Set-RowRange $ws 13 "G" @(5.780535,17.341605,0.7111490510422025,0.7111490510422023)
Set-RowRange $ws 13 "M" @(52.056859,156.170577,0.3786030350667928,0.3786030350667929,300.916495439565,2708.248458956085,0.2692431891094474,0.2692431891094474)
Set-RowRange $ws 14 "E" @(2,0.6666666666666666,0.1205213333333333,0.361564,0.01482711061006308,0.01482711061006307)
Set-RowRange $ws 14 "M" @(8.820647333333334,26.461942,0.06415146660411865,0.06415146660411865,1.063076177476444,9.567685597287999,0.0009511808911370346,0.0009511808911370344)
Set-RowRange $ws 15 "E" @(2,0.6666666666666666,0.1205213333333333,0.361564,0.01482711061006308,0.01482711061006307)
Set-RowRange $ws 15 "O" @(0.3979101621202897,0.3979101621202898,6.593907146913333,59.34516432222,0.005899857986625666,0.005899857986625666)
Set-RowRange $ws 16 "E" @(2,0.6666666666666666,0.1205213333333333,0.361564,0.01482711061006308,0.01482711061006307)
Set-RowRange $ws 16 "M" @(21.90816333333333,65.72449,0.1593353362087987,0.1593353362087987,2.640401055817778,23.76360950236,0.002362482654059446,0.002362482654059445)
Set-RowRange $ws 17 "E" @(2,0.6666666666666666,0.1205213333333333,0.361564,0.01482711061006308,0.01482711061006307)
Set-RowRange $ws 17 "M" @(52.056859,156.170577,0.3786030350667928,0.3786030350667929,6.273962055825333,56.46565850242799,0.005613589078240926,0.005613589078240926)
